$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style of the existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-21
$iValues = @(4, 9, 11, 8, 8, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1)
$jValues = @(6, 9, 12, 9, 9, 1, 5, 6, 6, 6, 5, 4, 6, 4, 6, 6, 3, 3, 4, 3)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
